# Add a "Greece" market sheet, cloned from the existing "Croatia" sheet
# (same layout/styles/merges), positioned right after it, with the
# Greece-specific Jira/ticket references filled in.

$wb = $excel.ActiveWorkbook

$croatia = $wb.Worksheets.Item("Croatia")

# Copy Croatia so the new sheet inherits identical formatting, column
# widths, merged cells and page setup; place the copy right after Croatia.
$croatia.Copy([System.Reflection.Missing]::Value, $croatia)

$greece = $wb.Worksheets.Item($wb.Worksheets.Count)
$greece.Name = "Greece"

# Fill in market-specific values (ticket id first, then market name, so
# the shared-string table is populated in the same order as the source
# edit: NGC-4119/T3202 before "Greece Market").
$greece.Range("B4").Value = "NGC-4119/T3202"
$greece.Range("B2").Value = "Greece Market"

# Restore Croatia's selection to a full-sheet selection (its previous
# "tabSelected" state is dropped once Greece becomes active).
$croatia.Activate()
$croatia.Cells.Select()

# Make Greece the active/selected sheet with B11 as the active cell.
$greece.Activate()
$greece.Range("B11").Select()
